# "changer Leo de groupe" - move Leo Nadeau from team in D2 to team in D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "Bryan Dubois<br/>Luka St-Hilaire<br/>Leo Nadeau" -> remove Leo Nadeau
$ws.Range("D2").Value = "Bryan Dubois<br/>Luka St-Hilaire"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Rows("2:2").RowHeight = 34

# D10: "Antoine Perreault<br/>Jayke Bédard<br/>Alexis Bergeron<br/>" -> add Leo Nadeau
$ws.Range("D10").Value = "Antoine Perreault<br/>Jayke Bédard<br/>Alexis Bergeron<br/>Leo Nadeau"
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Rows("10:10").RowHeight = 85

$excel.CutCopyMode = $false

$ws.Range("D11").Select() | Out-Null
